# 0.3.1.10b Updates with Permissions
#
# - Row 8 (Rateme): status goes from "In Progress" -> "Complete 0.3.1.9b"
#   (new shared string), picks up the green "Complete" fill/font used
#   elsewhere in column F, and gains a Line Complete value of 377.
# - Row 10 (Dice): status goes from "Incomplete" -> "In Progress", picking
#   up the yellow "In Progress" fill/font that row 8 used to have.
# - Row 4 Line Complete count bumped 340 -> 348.
# - Active selection moves from F4 to F8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats = -4122: copy only cell formatting (fill/font), not values.
$xlPasteFormats = -4122

# F10 is becoming "In Progress", which is the exact format F8 currently
# has, so grab that formatting first before F8's own format changes.
$ws.Range("F8").Copy() | Out-Null
$ws.Range("F10").PasteSpecial($xlPasteFormats) | Out-Null

# F8 is becoming a "Complete ..." status, so copy the green formatting
# already used by other "Complete" cells (e.g. F3) onto it.
$ws.Range("F3").Copy() | Out-Null
$ws.Range("F8").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Application.CutCopyMode = $false

# Update the cell values/text.
$ws.Range("F8").Value = "Complete 0.3.1.9b"
$ws.Range("G8").Value = 377
$ws.Range("F10").Value = "In Progress"
$ws.Range("G4").Value = 348

# Move the active selection to F8, matching the saved view state.
$ws.Range("F8").Select() | Out-Null
